$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "football youth compression pants"
$ws.Range("A2").Value = "knee guard"
$ws.Range("A3").Value = "knee sleeve basketball pair"
$ws.Range("A4").Value = "soccer tights for boys"
$ws.Range("A5").Value = "knee pads volleyball girls"
$ws.Range("A6").Value = "medias basketball"
$ws.Range("A7").Value = "boys yoga"
$ws.Range("A8").Value = "baseball chart"
$ws.Range("A9").Value = "basketball legs"
$ws.Range("A10").Value = "youth knee pad"
$ws.Range("A11").Value = "softball band"
$ws.Range("A12").Value = "need pads volleyball"
$ws.Range("A13").Value = "mens yoga pants capri"
$ws.Range("A14").Value = "compression shorts knee length"
$ws.Range("A15").Value = "men thigh compression pants"
$ws.Range("A16").Value = "basketball leg compression"
$ws.Range("A17").Value = "boys youth leggings"
$ws.Range("A18").Value = "fabric softball"
$ws.Range("A19").Value = "bump pad"
$ws.Range("A20").Value = "men knee pads for work"
$ws.Range("A21").Value = "black spandex pants men"
$ws.Range("A22").Value = "youth volleyball kneepads"
$ws.Range("A23").Value = "compression tight men"
$ws.Range("A24").Value = "black compression tights men"
$ws.Range("A25").Value = "basketball sleeve with pad"
$ws.Range("A26").Value = "compression padded knee sleeve"
$ws.Range("A27").Value = "compression pants black"
$ws.Range("A28").Value = "men capris"
$ws.Range("A29").Value = "calf pad"
$ws.Range("A30").Value = "knee pads volleyball girls youth"
$ws.Range("A31").Value = "men yoga pants"
$ws.Range("A32").Value = "slim knee pads"
$ws.Range("A33").Value = "hockey leg pads"
$ws.Range("A34").Value = "soccer leggings men"
$ws.Range("A35").Value = "knee pads volleyball large"
$ws.Range("A36").Value = "softball fabric"
$ws.Range("A37").Value = "knee pad thick"
$ws.Range("A38").Value = "yoga for knees"
$ws.Range("A39").Value = "knee pads for boys"
$ws.Range("A40").Value = "knee sleeve honeycomb"
$ws.Range("A41").Value = "black mens leggings"
$ws.Range("A42").Value = "used softballs"
$ws.Range("A43").Value = "youth tights for sports"
$ws.Range("A44").Value = "boy compression leggings"
$ws.Range("A45").Value = "cold compression pants"
$ws.Range("A46").Value = "youth knee guards"
$ws.Range("A47").Value = "paintball knee"
$ws.Range("A48").Value = "knee pad construction"
$ws.Range("A49").Value = "basketball knee sleeve youth"
$ws.Range("A50").Value = "mens wrestling shorts"
$ws.Range("A51").Value = "recovery compression tights men"
$ws.Range("A52").Value = "womens sliding shorts softball"
$ws.Range("A53").Value = "gym kneeling pad"
$ws.Range("A54").Value = "gym leggings for men"
$ws.Range("A55").Value = "durable pants"
$ws.Range("A56").Value = "adult leggings"
$ws.Range("A57").Value = "work wear knee pads"
$ws.Range("A58").Value = "black lacrosse shorts"
$ws.Range("A59").Value = "hex skin padding"
$ws.Range("A60").Value = "sport leggings boys"
$ws.Range("A61").Value = "need pads for construction"
$ws.Range("A62").Value = "football girdle youth"
$ws.Range("A63").Value = "5 inch foam basketball"
$ws.Range("A64").Value = "mens skin tight leggings"
$ws.Range("A65").Value = "silicon knee pads"
$ws.Range("A66").Value = "extra small baseball pants"
$ws.Range("A67").Value = "football pads for men"
$ws.Range("A68").Value = "mens kneepads"
$ws.Range("A69").Value = "thigh protector men"
$ws.Range("A70").Value = "long volleyball knee pads"
$ws.Range("A71").Value = "paintball pants youth"
$ws.Range("A72").Value = "gym leggings men"
$ws.Range("A73").Value = "compression running capris"
$ws.Range("A74").Value = "youth running pants boys"
$ws.Range("A75").Value = "knee sleeve baseball"
$ws.Range("A76").Value = "football pads youth"
$ws.Range("A77").Value = "volleyball spandex shorts"
$ws.Range("A78").Value = "work in baseball"
$ws.Range("A79").Value = "thick yoga knee pad"
$ws.Range("A80").Value = "running tights mens"
$ws.Range("A81").Value = "basketballs 28 5"
$ws.Range("A82").Value = "softball pants"
$ws.Range("A83").Value = "youth volleyball spandex"
$ws.Range("A84").Value = "work pants knee pads"
$ws.Range("A85").Value = "size small baseball pants"
$ws.Range("A86").Value = "boys sport tights"
$ws.Range("A87").Value = "mens medium tall athletic pants"
$ws.Range("A88").Value = "compression sleeve knee youth"
$ws.Range("A89").Value = "compression knee sleeves for basketball"
$ws.Range("A90").Value = "mens knee pads for work"
$ws.Range("A91").Value = "adult knee pads for work"
$ws.Range("A92").Value = "leggings cycling"
$ws.Range("A93").Value = "male athletic tights"
$ws.Range("A94").Value = "knee compression sleeve basketball"
$ws.Range("A95").Value = "compressions knee"
$ws.Range("A96").Value = "youth running tights"
$ws.Range("A97").Value = "knee chart"
$ws.Range("A98").Value = "long softball pants"
$ws.Range("A99").Value = "compression pants men soccer"
$ws.Range("A100").Value = "limber support"
